# Apply the edits described by the diff:
# 1. Change a few "missing" (empty) vs filled values in column F for several rows
# 2. Remove the "RM 232" row and the "SC 92" row entirely (rows shift up)
# 3. Dimension will naturally update from A1:F35 to A1:F33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: toggle a handful of F-column values on the current (pre-delete) layout ---

# F3 (RM 8): was empty -> now 17.64
$ws.Range("F3").Value = 17.64

# F5 (RM 14): was 17.66 -> now empty
$ws.Range("F5").ClearContents()

# F21 (RM 135): was empty -> now 16.58
$ws.Range("F21").Value = 16.58

# F23 (RM 140): was 16.48 -> now empty
$ws.Range("F23").ClearContents()

# --- Step 2: delete the "RM 232" row (row 26) entirely ---
$ws.Rows("26").Delete()

# After that deletion, the row that used to be "SC 92" (originally row 28)
# is now row 27. Delete it too.
$ws.Rows("27").Delete()

# --- Step 3: the row that is now "SC 193" (originally row 34, now row 32)
# previously had an empty F value; it now gets a value of 17.39
$ws.Range("F32").Value = 17.39
